$wb = $excel.ActiveWorkbook

# --- Sheet 1 ("Purchase 22-23"): remove the "Sanyo and Sanyo" entry (row 10) and
# the two "Asha Enterprises" entries (rows 14-15), pulling the rows below up so
# the sheet keeps its usual data/blank-row rhythm. Deleting whole rows via COM
# also purges the now-unused shared strings automatically (matches the diff's
# sharedStrings.xml shrink from 45 -> 43 unique strings).
$ws1 = $wb.Worksheets.Item(1)

# Remove the "Sanyo and Sanyo" row (10) together with its blank spacer row (11).
$ws1.Rows("10:11").Delete()

# Remove both "Asha Enterprises" rows (now at 12-13 after the shift above) plus
# the following blank spacer row (14).
$ws1.Rows("12:14").Delete()

# The rows that slid up keep their original "Sr. No" values from the diff
# (the source workbook already had a small numbering quirk - row 12 repeats
# "6" and row 14 keeps "7" - so we restore those exact numbers here).
$ws1.Cells.Item(12, 1).Value = 6
$ws1.Cells.Item(14, 1).Value = 7

# --- View state: the edit also flips which sheet/tab is active and updates the
# remembered selections on both sheets.
$ws2 = $wb.Worksheets.Item(2)

$ws1.Activate()
$ws1.Range("A1:G14").Select()

$ws2.Activate()
$ws2.Range("E29").Select()
